# Re-resolve the workbook/sheet from $excel explicitly (the pre-seeded
# $wb binding is unreliable in this host), then grab the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case rows (Entrada / Saida) appended below the existing header row.
$data = @(
    @("(1000;5)", "SL 920; 80; 80; -27,95"),
    @("(1500;3)", "1380; 120;120;811,23"),
    @("2500;5", "2275; 225; 200; 947,95"),
    @("1200;2", "1104; 96;96;379,18"),
    @("3600;3", "3146,16; 396;288;54,84"),
    @("1700;1", "1492,16;153;136;54,84"),
    @("1800;2", "1638;162;144;0"),
    @("4000;5", "3506,09;440;320;53,10"),
    @("7000;3", "5650,76;621,04;560;728,20"),
    @("1780;2", "891,60;160,20;142,40;728,20")
)

$row = 3
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Leave the cursor where the author ended up after typing the last row.
[void]$ws.Range("B13").Select()

# Columns widened (e.g. after content-driven auto-fit) for the four used columns.
$ws.Columns.Item(1).ColumnWidth = 54.666666666666664
$ws.Columns.Item(2).ColumnWidth = 56.166666666666664
$ws.Columns.Item(3).ColumnWidth = 34.333333333333336
$ws.Columns.Item(4).ColumnWidth = 29.0

# Best-effort: persist the tab ratio seen in the target workbook view.
try {
    $excel.Windows.Item(1).TabRatio = 0.993
} catch {
}
$null = $null
